# Fixing geopoint in shared_table model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Insert a new column before column C (schema.elementType shifts right,
# making room for the new "schema.name" column).
$ws.Columns("C").Insert()

# New header cell for the inserted column.
$ws.Range("C1").Value = "schema.name"

# The geopoint sub-property headers now describe the property's "type",
# so append ".type" to each (columns shifted E..H after the insert).
$ws.Range("E1").Value = "schema.properties.latitude.type"
$ws.Range("F1").Value = "schema.properties.longitude.type"
$ws.Range("G1").Value = "schema.properties.altitude.type"
$ws.Range("H1").Value = "schema.properties.accuracy.type"

# Row 4 (refrigerator_location): schema.name should mirror the
# schema.elementType value ("geopoint"), same as column D after the shift.
$ws.Range("C4").Value = $ws.Range("D4").Text
